$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a literal text value to a cell without letting Excel
# auto-convert numeric-looking strings (e.g. "574.60") into numbers, and
# without leaving a permanent style change behind on the cell.
function Set-TextValue($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# --- Update "Price" (column D) and "Volume(1h)" (column E) values ---
# Row 2
Set-TextValue $ws "D2" '66.082.07'
$ws.Range("E2").Value = '  +2.66%  '
# Row 3
Set-TextValue $ws "D3" '2.963.85'
$ws.Range("E3").Value = '  -0.24%  '
# Row 4
Set-TextValue $ws "D4" '0.999'
$ws.Range("E4").Value = '  -0.02%  '
# Row 5
Set-TextValue $ws "D5" '574.60'
$ws.Range("E5").Value = '  -0.78%  '
# Row 6
Set-TextValue $ws "D6" '161.26'
$ws.Range("E6").Value = '  +5.60%  '
# Row 7
$ws.Range("E7").Value = '  -0.02%  '
# Row 8
$ws.Range("E8").Value = '  +0.93%  '
# Row 9
Set-TextValue $ws "D9" '2.960.14'
$ws.Range("E9").Value = '  -0.30%  '
# Row 10
$ws.Range("E10").Value = '  -5.03%  '
# Row 11
$ws.Range("E11").Value = '  -0.34%  '
# Row 12
Set-TextValue $ws "D12" '0.455'
$ws.Range("E12").Value = '  +2.26%  '
# Row 13
Set-TextValue $ws "D13" '0.0000246'
$ws.Range("E13").Value = '  +1.92%  '
# Row 14
Set-TextValue $ws "D14" '34.30'
$ws.Range("E14").Value = '  -0.41%  '
# Row 15
$ws.Range("E15").Value = '  -0.74%  '
# Row 16
Set-TextValue $ws "D16" '66.070.60'
$ws.Range("E16").Value = '  +2.77%  '
# Row 17
Set-TextValue $ws "D17" '3.455.08'
$ws.Range("E17").Value = '  -0.20%  '
# Row 18
Set-TextValue $ws "D18" '6.90'
$ws.Range("E18").Value = '  +0.20%  '
# Row 19
Set-TextValue $ws "D19" '2.964.25'
$ws.Range("E19").Value = '  -0.12%  '
# Row 20
Set-TextValue $ws "D20" '448.13'
$ws.Range("E20").Value = '  +0.74%  '
# Row 21
Set-TextValue $ws "D21" '13.74'
$ws.Range("E21").Value = '  +1.16%  '
# Row 22
Set-TextValue $ws "D22" '0.675'
$ws.Range("E22").Value = '  -0.19%  '
# Row 23
Set-TextValue $ws "D23" '7.16'
$ws.Range("E23").Value = '  -0.51%  '
# Row 24
Set-TextValue $ws "D24" '81.91'
$ws.Range("E24").Value = '  +1.66%  '
# Row 25
$ws.Range("E25").Value = '  +1.40%  '
# Row 26
Set-TextValue $ws "D26" '12.17'
$ws.Range("E26").Value = '  -1.05%  '
# Row 27
$ws.Range("E27").Value = '  -0.12%  '
# Row 28
Set-TextValue $ws "D28" '10.01'
$ws.Range("E28").Value = '  -8.80%  '
# Row 29
$ws.Range("E29").Value = '  +4.82%  '
# Row 30
Set-TextValue $ws "D30" '2.39'
$ws.Range("E30").Value = '  +8.74%  '
# Row 31
$ws.Range("E31").Value = '  +0.82%  '
# Row 32
Set-TextValue $ws "D32" '0.0₃0991'
$ws.Range("E32").Value = '  -9.16%  '
# Row 33
Set-TextValue $ws "D33" '27.20'
$ws.Range("E33").Value = '  +2.66%  '
# Row 34
$ws.Range("E34").Value = '  -0.42%  '
# Row 35
$ws.Range("E35").Value = '  +0.05%  '
# Row 36
$ws.Range("E36").Value = '  +0.25%  '
# Row 37
Set-TextValue $ws "D37" '5.72'
$ws.Range("E37").Value = '  +1.45%  '
# Row 38
Set-TextValue $ws "D38" '49.35'
$ws.Range("E38").Value = '  +0.86%  '
# Row 39
$ws.Range("E39").Value = '  -4.27%  '
# Row 40
Set-TextValue $ws "D40" '43.40'
$ws.Range("E40").Value = '  -1.05%  '
# Row 41
Set-TextValue $ws "D41" '0.299'
$ws.Range("E41").Value = '  +2.44%  '
# Row 42
Set-TextValue $ws "D42" '2.83'
$ws.Range("E42").Value = '  -7.79%  '
# Row 43
$ws.Range("E43").Value = '  +0.10%  '
# Row 44
Set-TextValue $ws "D44" '8.37'
$ws.Range("E44").Value = '  +0.45%  '
# Row 45
Set-TextValue $ws "D45" '385.04'
$ws.Range("E45").Value = '  -0.49%  '
# Row 46
Set-TextValue $ws "D46" '0.0354'
$ws.Range("E46").Value = '  +1.50%  '
# Row 47
Set-TextValue $ws "D47" '2.712.16'
$ws.Range("E47").Value = '  -1.87%  '
# Row 48
Set-TextValue $ws "D48" '131.05'
$ws.Range("E48").Value = '  -2.62%  '
# Row 49
$ws.Range("E49").Value = '  +0.04%  '
# Row 50
$ws.Range("E50").Value = '  +0.54%  '
# Row 51
$ws.Range("E51").Value = '  +0.13%  '
